# This script regenerates the "K" (strikeouts) column (column G) of the save-data
# sheet so that it reflects the corrected per-game values (previously this column
# held stale "Strike#" figures). All other columns are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value for column G
$kValues = @{
    2 = 1
    3 = 0
    4 = 1
    5 = 1
    6 = 2
    8 = 2
    9 = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 2
    17 = 3
    18 = 1
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 3
    25 = 1
    26 = 1
    27 = 1
    28 = 3
    29 = 2
    30 = 2
    31 = 0
    32 = 1
    33 = 1
    34 = 3
    35 = 0
    36 = 1
    37 = 1
    38 = 0
    39 = 0
    40 = 0
    41 = 1
    42 = 2
    43 = 1
    44 = 2
    45 = 0
    46 = 2
    47 = 1
    48 = 0
    49 = 0
    50 = 1
    51 = 3
    52 = 1
    53 = 4
    54 = 1
    55 = 1
    56 = 0
    57 = 1
    58 = 2
    61 = 1
    62 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}
